$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.397.98'
$ws.Range("E2").Value = '  -0.42%  '
$ws.Range("D3").Value = '1.846.81'
$ws.Range("E3").Value = '  -0.25%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9985'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.07'
$ws.Range("E5").Value = '  -1.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6326'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07553'
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2964'
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.59'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07736'
$ws.Range("E11").Value = '  +0.72%  '
$ws.Range("D12").Value = '1.850.27'
$ws.Range("E12").Value = '  -0.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.999'
$ws.Range("E13").Value = '  -0.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6845'
$ws.Range("E14").Value = '  -0.26%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001002'
$ws.Range("E15").Value = '  +1.73%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.05'
$ws.Range("E16").Value = '  -0.96%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.183'
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").Value = '29.424.36'
$ws.Range("E18").Value = '  -0.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '229.58'
$ws.Range("E19").Value = '  -2.35%  '
$ws.Range("E20").Value = '  -0.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9994'
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.573'
$ws.Range("E22").Value = '  -0.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.961'
$ws.Range("E23").Value = '  -0.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.11'
$ws.Range("E25").Value = '  +0.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1401'
$ws.Range("E26").Value = '  +1.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.380'
$ws.Range("E27").Value = '  -0.81%  '
$ws.Range("E28").Value = '  -0.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.463'
$ws.Range("E29").Value = '  -1.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05727'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.251'
$ws.Range("E31").Value = '  -1.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.130'
$ws.Range("E32").Value = '  +0.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.038'
$ws.Range("E33").Value = '  -0.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.855'
$ws.Range("E34").Value = '  -2.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.156'
$ws.Range("E35").Value = '  -1.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7184'
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").Value = '1.251.56'
$ws.Range("E38").Value = '  +1.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01808'
$ws.Range("E39").Value = '  +1.85%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.781'
$ws.Range("E40").Value = '  -0.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.196'
$ws.Range("E41").Value = '  +0.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9061'
$ws.Range("E42").Value = '  -0.94%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").Value = '2.008.04'
$ws.Range("E44").Value = '  -1.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.77'
$ws.Range("E45").Value = '  -0.27%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '66.44'
$ws.Range("E46").Value = '  -1.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.073'
$ws.Range("E47").Value = '  -3.62%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000118'
$ws.Range("E48").Value = '  +0.24%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.193'
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4027'
$ws.Range("E50").Value = '  -0.30%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.711'
$ws.Range("E51").Value = '  +1.19%  '
